$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Range("B2").Value = "Bitcoin"
$ws.Range("C2").Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
Set-TextCell $ws.Range("D2") "61.598.26"
$ws.Range("E2").Value = "  -3.84%  "

$ws.Range("B3").Value = "Ethereum"
$ws.Range("C3").Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
Set-TextCell $ws.Range("D3") "2.975.22"
$ws.Range("E3").Value = "  -5.09%  "

$ws.Range("B4").Value = "TetherUSD"
$ws.Range("C4").Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
Set-TextCell $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
Set-TextCell $ws.Range("D5") "542.43"
$ws.Range("E5").Value = "  -4.90%  "

$ws.Range("B6").Value = "Solana"
$ws.Range("C6").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextCell $ws.Range("D6") "152.19"
$ws.Range("E6").Value = "  -5.85%  "

$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextCell $ws.Range("D7") "0.999"
$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextCell $ws.Range("D8") "0.575"
$ws.Range("E8").Value = "  +0.70%  "

$ws.Range("B9").Value = "LidoStakedEther"
$ws.Range("C9").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
Set-TextCell $ws.Range("D9") "2.986.71"
$ws.Range("E9").Value = "  -5.13%  "

$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextCell $ws.Range("D10") "0.114"
$ws.Range("E10").Value = "  -2.01%  "

$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell $ws.Range("D11") "6.15"
$ws.Range("E11").Value = "  -6.61%  "

$ws.Range("B12").Value = "Cardano"
$ws.Range("C12").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextCell $ws.Range("D12") "0.371"
$ws.Range("E12").Value = "  -3.39%  "

$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextCell $ws.Range("D13") "3.501.13"
$ws.Range("E13").Value = "  -5.09%  "

$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextCell $ws.Range("D14") "0.125"
$ws.Range("E14").Value = "  -2.34%  "

$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextCell $ws.Range("D15") "61.664.99"
$ws.Range("E15").Value = "  -3.94%  "

$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextCell $ws.Range("D16") "23.74"
$ws.Range("E16").Value = "  -5.12%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell $ws.Range("D17") "2.987.67"
$ws.Range("E17").Value = "  -5.03%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextCell $ws.Range("D18") "0.0000147"
$ws.Range("E18").Value = "  -4.15%  "

$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell $ws.Range("D19") "5.18"
$ws.Range("E19").Value = "  -1.09%  "

$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell $ws.Range("D20") "12.04"
$ws.Range("E20").Value = "  -4.10%  "

$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextCell $ws.Range("D21") "381.72"
$ws.Range("E21").Value = "  -4.91%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextCell $ws.Range("D22") "6.70"
$ws.Range("E22").Value = "  -6.14%  "

$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextCell $ws.Range("D23") "1.00"
$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextCell $ws.Range("D24") "5.66"
$ws.Range("E24").Value = "  -3.39%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell $ws.Range("D25") "65.71"
$ws.Range("E25").Value = "  -3.00%  "

$ws.Range("B26").Value = "Polygon"
$ws.Range("C26").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextCell $ws.Range("D26") "0.471"
$ws.Range("E26").Value = "  -2.46%  "

$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextCell $ws.Range("D27") "3.102.58"
$ws.Range("E27").Value = "  -5.21%  "

$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell $ws.Range("D28") "0.190"
$ws.Range("E28").Value = "  -1.41%  "

$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextCell $ws.Range("D29") "0.0₃0947"
$ws.Range("E29").Value = "  -5.58%  "

$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextCell $ws.Range("D30") "0.998"
$ws.Range("E30").Value = "  +0.66%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell $ws.Range("D31") "8.26"
$ws.Range("E31").Value = "  -6.06%  "

$ws.Range("B32").Value = "USDe"
$ws.Range("C32").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextCell $ws.Range("D32") "0.999"
$ws.Range("E32").Value = "  +0.01%  "

$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell $ws.Range("D33") "1.72"
$ws.Range("E33").Value = "  -4.70%  "

$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell $ws.Range("D34") "20.56"
$ws.Range("E34").Value = "  -2.81%  "

$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell $ws.Range("D35") "160.72"
$ws.Range("E35").Value = "  +0.60%  "

$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell $ws.Range("D36") "4.68"
$ws.Range("E36").Value = "  -2.65%  "

$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell $ws.Range("D37") "5.95"
$ws.Range("E37").Value = "  -5.00%  "

$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell $ws.Range("D38") "1.08"
$ws.Range("E38").Value = "  -2.92%  "

$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell $ws.Range("D39") "1.27"
$ws.Range("E39").Value = "  -5.04%  "

$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell $ws.Range("D40") "1.56"
$ws.Range("E40").Value = "  -6.10%  "

$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell $ws.Range("D41") "3.94"
$ws.Range("E41").Value = "  -3.14%  "

$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell $ws.Range("D42") "2.417.88"
$ws.Range("E42").Value = "  -9.38%  "

$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell $ws.Range("D43") "37.47"
$ws.Range("E43").Value = "  -2.19%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell $ws.Range("D44") "22.18"
$ws.Range("E44").Value = "  -6.27%  "

$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextCell $ws.Range("D45") "0.671"
$ws.Range("E45").Value = "  -2.72%  "

$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell $ws.Range("D46") "0.0595"
$ws.Range("E46").Value = "  -2.99%  "

$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell $ws.Range("D47") "5.18"
$ws.Range("E47").Value = "  -4.73%  "

$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell $ws.Range("D48") "0.0248"
$ws.Range("E48").Value = "  -2.75%  "

$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextCell $ws.Range("D49") "0.997"
$ws.Range("E49").Value = "  -0.02%  "

$ws.Range("B50").Value = "Bittensor"
$ws.Range("C50").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextCell $ws.Range("D50") "270.54"
$ws.Range("E50").Value = "  -5.92%  "

$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell $ws.Range("D51") "0.0956"
$ws.Range("E51").Value = "  -2.18%  "
